$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.623.75'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.493.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.13%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.30%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.69%  '

$ws.Range("E7").Value = '  -0.66%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("E9").Value = '  -0.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.84%  '

$ws.Range("E11").Value = '  -0.63%  '

$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.73%  '

$ws.Range("E14").Value = '  +0.58%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.886.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.28%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.493.49'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.847'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.488.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.63'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.74%  '

$ws.Range("E21").Value = '  +0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.74'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +14.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '247.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.64%  '

$ws.Range("E25").Value = '  -1.13%  '

$ws.Range("E26").Value = '  -0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.57%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.31'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.76%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.139'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.20%  '

$ws.Range("E33").Value = '  +2.10%  '

$ws.Range("E34").Value = '  -0.65%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0787'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.42%  '

$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.63%  '

$ws.Range("E38").Value = '  -0.52%  '

$ws.Range("E39").Value = '  -1.38%  '

$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '119.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.71%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0298'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.994.26'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.36%  '

$ws.Range("E46").Value = '  +1.48%  '

$ws.Range("E47").Value = '  -2.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.80'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.24%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.57%  '

$ws.Range("E50").Value = '  -0.84%  '

$ws.Range("E51").Value = '  +3.26%  '
